$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 79, shifting existing rows 79..162 down to 80..163
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new record
$ws.Cells.Item(79, 1).Value = 7
$ws.Cells.Item(79, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(79, 3).Value = "Ñuble"
$ws.Cells.Item(79, 4).Value = 44494
$ws.Cells.Item(79, 5).Value = 16
$ws.Cells.Item(79, 6).Value = 100112043
$ws.Cells.Item(79, 7).Value = "Pepino ensalada"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 120
$ws.Cells.Item(79, 11).Value = 13000
$ws.Cells.Item(79, 12).Value = 14000
$ws.Cells.Item(79, 13).Value = 13500
$ws.Cells.Item(79, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(79, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(79, 16).Value = 225
$ws.Cells.Item(79, 17).Value = 60
$ws.Cells.Item(79, 18).Value = "Hortaliza"
